# Updates the crypto price/volume table on Sheet1 with refreshed figures.
# D-column price cells are forced to Text (NumberFormat "@") before the
# assignment and then ClearFormats() is used to drop the explicit format
# again, because several of the new price strings (e.g. "38.90", "0.101")
# are otherwise auto-coerced by Excel into numeric values, which would
# both change the stored type and silently drop significant trailing
# zeros / digits. Row 48/49 additionally swap Coin/Link/Price content
# (Stellar now ranked above Mantle).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '63.891.62'
$c.ClearFormats()
$ws.Range('E2').Value = '  -0.61%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.747.49'
$c.ClearFormats()
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('E4').Value = '  +0.08%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '572.88'
$c.ClearFormats()
$ws.Range('E5').Value = '  -1.41%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '157.03'
$c.ClearFormats()
$ws.Range('E6').Value = '  +1.23%  '
$ws.Range('E7').Value = '  -0.06%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.602'
$c.ClearFormats()
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('E9').Value = '  -2.98%  '
$ws.Range('E10').Value = '  +0.08%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.381'
$c.ClearFormats()
$ws.Range('E11').Value = '  -2.06%  '
$ws.Range('E12').Value = '  -18.84%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '3.234.16'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.85%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '26.42'
$c.ClearFormats()
$ws.Range('E14').Value = '  -0.63%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '63.573.41'
$c.ClearFormats()
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('E16').Value = '  -2.06%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '2.752.36'
$c.ClearFormats()
$ws.Range('E17').Value = '  -0.69%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '12.15'
$c.ClearFormats()
$ws.Range('E18').Value = '  +1.11%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '4.79'
$c.ClearFormats()
$ws.Range('E19').Value = '  -2.16%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '353.93'
$c.ClearFormats()
$ws.Range('E20').Value = '  -2.17%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.73'
$c.ClearFormats()
$ws.Range('E21').Value = '  -4.35%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.536'
$c.ClearFormats()
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('E23').Value = '  -0.47%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '65.17'
$c.ClearFormats()
$ws.Range('E24').Value = '  -2.17%  '
$ws.Range('E25').Value = '  -1.77%  '
$ws.Range('E26').Value = '  -0.12%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '8.38'
$c.ClearFormats()
$ws.Range('E27').Value = '  -1.49%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.0₃0899'
$c.ClearFormats()
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('E29').Value = '  -4.07%  '
$ws.Range('E30').Value = '  -2.52%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '169.47'
$c.ClearFormats()
$ws.Range('E31').Value = '  -1.23%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.20'
$c.ClearFormats()
$ws.Range('E32').Value = '  -7.44%  '
$ws.Range('E33').Value = '  -2.06%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('E36').Value = '  -0.24%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.78'
$c.ClearFormats()
$ws.Range('E37').Value = '  -2.42%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.976'
$c.ClearFormats()
$ws.Range('E38').Value = '  -3.83%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '6.13'
$c.ClearFormats()
$ws.Range('E39').Value = '  +4.95%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '4.12'
$c.ClearFormats()
$ws.Range('E40').Value = '  -2.88%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '325.35'
$c.ClearFormats()
$ws.Range('E41').Value = '  -6.20%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '38.90'
$c.ClearFormats()
$ws.Range('E42').Value = '  -0.89%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '21.26'
$c.ClearFormats()
$ws.Range('E43').Value = '  -2.91%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0585'
$c.ClearFormats()
$ws.Range('E44').Value = '  -0.97%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '21.26'
$c.ClearFormats()
$ws.Range('E45').Value = '  -3.16%  '
$ws.Range('E46').Value = '  -1.46%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '134.82'
$c.ClearFormats()
$ws.Range('E47').Value = '  -2.14%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.101'
$c.ClearFormats()
$ws.Range('E48').Value = '  -0.71%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.622'
$c.ClearFormats()
$ws.Range('E49').Value = '  -4.47%  '
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('E51').Value = '  +0.42%  '

Write-Host "Applied all changes"